# Auto-generated edit script applying numeric updates from the commit diff
$wb = $excel.ActiveWorkbook

$sheets = @{
    "ALC" = $wb.Worksheets.Item("ALC")
    "ARM" = $wb.Worksheets.Item("ARM")
    "BSM" = $wb.Worksheets.Item("BSM")
    "CRP" = $wb.Worksheets.Item("CRP")
    "CUL" = $wb.Worksheets.Item("CUL")
    "GSM" = $wb.Worksheets.Item("GSM")
    "LTW" = $wb.Worksheets.Item("LTW")
    "WVR" = $wb.Worksheets.Item("WVR")
}

$sheets["ALC"].Range("H43").Value = 563.05884
$sheets["ALC"].Range("I43").Value = 581.5833
$sheets["ALC"].Range("J43").Value = 518.6
$sheets["ALC"].Range("K43").Value = 581.5833
$sheets["ALC"].Range("L43").Value = 518.6
$sheets["ALC"].Range("M43").Value = -512.5833
$sheets["ALC"].Range("N43").Value = -656.6
$sheets["ALC"].Range("H64").Value = 4061.8572
$sheets["ALC"].Range("I64").Value = 3450
$sheets["ALC"].Range("K64").Value = 3450
$sheets["ALC"].Range("M64").Value = -3202
$sheets["ALC"].Range("H67").Value = 4061.8572
$sheets["ALC"].Range("I67").Value = 3450
$sheets["ALC"].Range("K67").Value = 3450
$sheets["ALC"].Range("M67").Value = -2592
$sheets["ALC"].Range("H74").Value = 5149.8335
$sheets["ALC"].Range("I74").Value = 4533.1665
$sheets["ALC"].Range("J74").Value = 5766.5
$sheets["ALC"].Range("K74").Value = 4533.1665
$sheets["ALC"].Range("L74").Value = 5766.5
$sheets["ALC"].Range("M74").Value = -3597.1665
$sheets["ALC"].Range("N74").Value = -7638.5
$sheets["ALC"].Range("H77").Value = 5149.8335
$sheets["ALC"].Range("I77").Value = 4533.1665
$sheets["ALC"].Range("J77").Value = 5766.5
$sheets["ALC"].Range("K77").Value = 22665.8325
$sheets["ALC"].Range("L77").Value = 28832.5
$sheets["ALC"].Range("M77").Value = -17985.8325
$sheets["ALC"].Range("N77").Value = -38192.5
$sheets["ALC"].Range("H100").Value = 1505.25
$sheets["ALC"].Range("I100").Value = 1469.8572
$sheets["ALC"].Range("J100").Value = 1753
$sheets["ALC"].Range("K100").Value = 1469.8572
$sheets["ALC"].Range("L100").Value = 1753
$sheets["ALC"].Range("M100").Value = -928.8571999999999
$sheets["ALC"].Range("N100").Value = -2835
$sheets["ARM"].Range("H2").Value = 1866.591
$sheets["ARM"].Range("I2").Value = 2014.0667
$sheets["ARM"].Range("J2").Value = 1550.5714
$sheets["ARM"].Range("K2").Value = 2014.0667
$sheets["ARM"].Range("L2").Value = 1550.5714
$sheets["ARM"].Range("M2").Value = -1901.0667
$sheets["ARM"].Range("N2").Value = -1776.5714
$sheets["ARM"].Range("H45").Value = 1913.1578
$sheets["ARM"].Range("I45").Value = 1843.75
$sheets["ARM"].Range("J45").Value = 2283.3333
$sheets["ARM"].Range("K45").Value = 1843.75
$sheets["ARM"].Range("L45").Value = 2283.3333
$sheets["ARM"].Range("M45").Value = -1466.75
$sheets["ARM"].Range("N45").Value = -3037.3333
$sheets["ARM"].Range("H63").Value = 3359.875
$sheets["ARM"].Range("I63").Value = 2669.6667
$sheets["ARM"].Range("J63").Value = 3774
$sheets["ARM"].Range("K63").Value = 2669.6667
$sheets["ARM"].Range("L63").Value = 3774
$sheets["ARM"].Range("M63").Value = -1983.6667
$sheets["ARM"].Range("N63").Value = -5146
$sheets["ARM"].Range("H66").Value = 3359.875
$sheets["ARM"].Range("I66").Value = 2669.6667
$sheets["ARM"].Range("J66").Value = 3774
$sheets["ARM"].Range("K66").Value = 13348.3335
$sheets["ARM"].Range("L66").Value = 18870
$sheets["ARM"].Range("M66").Value = -9916.333500000001
$sheets["ARM"].Range("N66").Value = -25734
$sheets["ARM"].Range("H88").Value = 6171.7144
$sheets["ARM"].Range("I88").Value = 11755.2
$sheets["ARM"].Range("J88").Value = 3069.7778
$sheets["ARM"].Range("K88").Value = 11755.2
$sheets["ARM"].Range("L88").Value = 3069.7778
$sheets["ARM"].Range("M88").Value = -11349.2
$sheets["ARM"].Range("N88").Value = -3881.7778
$sheets["ARM"].Range("H91").Value = 6171.7144
$sheets["ARM"].Range("I91").Value = 11755.2
$sheets["ARM"].Range("J91").Value = 3069.7778
$sheets["ARM"].Range("K91").Value = 11755.2
$sheets["ARM"].Range("L91").Value = 3069.7778
$sheets["ARM"].Range("M91").Value = -10351.2
$sheets["ARM"].Range("N91").Value = -5877.7778
$sheets["ARM"].Range("H102").Value = 2018.3572
$sheets["ARM"].Range("I102").Value = 1621.4166
$sheets["ARM"].Range("J102").Value = 4400
$sheets["ARM"].Range("K102").Value = 1621.4166
$sheets["ARM"].Range("L102").Value = 4400
$sheets["ARM"].Range("M102").Value = 0.5833999999999833
$sheets["ARM"].Range("N102").Value = -7644
$sheets["ARM"].Range("H110").Value = 1420
$sheets["ARM"].Range("I110").Value = 1355.5555
$sheets["ARM"].Range("J110").Value = 2000
$sheets["ARM"].Range("K110").Value = 1355.5555
$sheets["ARM"].Range("L110").Value = 2000
$sheets["ARM"].Range("M110").Value = 689.4445000000001
$sheets["ARM"].Range("N110").Value = -6090
$sheets["ARM"].Range("H116").Value = 1866.591
$sheets["ARM"].Range("I116").Value = 2014.0667
$sheets["ARM"].Range("J116").Value = 1550.5714
$sheets["ARM"].Range("K116").Value = 2014.0667
$sheets["ARM"].Range("L116").Value = 1550.5714
$sheets["ARM"].Range("M116").Value = 279.9332999999999
$sheets["ARM"].Range("N116").Value = -6138.5714
$sheets["BSM"].Range("H3").Value = 1866.591
$sheets["BSM"].Range("I3").Value = 2014.0667
$sheets["BSM"].Range("J3").Value = 1550.5714
$sheets["BSM"].Range("K3").Value = 2014.0667
$sheets["BSM"].Range("L3").Value = 1550.5714
$sheets["BSM"].Range("M3").Value = -1900.0667
$sheets["BSM"].Range("N3").Value = -1778.5714
$sheets["BSM"].Range("H86").Value = 1854.4222
$sheets["BSM"].Range("I86").Value = 1844.5
$sheets["BSM"].Range("K86").Value = 1844.5
$sheets["BSM"].Range("M86").Value = -721.5
$sheets["BSM"].Range("H89").Value = 1854.4222
$sheets["BSM"].Range("I89").Value = 1844.5
$sheets["BSM"].Range("K89").Value = 9222.5
$sheets["BSM"].Range("M89").Value = -3606.5
$sheets["BSM"].Range("H94").Value = 1007.5
$sheets["BSM"].Range("I94").Value = 1129
$sheets["BSM"].Range("K94").Value = 1129
$sheets["BSM"].Range("M94").Value = -678
$sheets["BSM"].Range("H105").Value = 8312.333000000001
$sheets["BSM"].Range("I105").Value = 8466.666999999999
$sheets["BSM"].Range("J105").Value = 8003.6665
$sheets["BSM"].Range("K105").Value = 8466.666999999999
$sheets["BSM"].Range("L105").Value = 8003.6665
$sheets["BSM"].Range("M105").Value = -6719.666999999999
$sheets["BSM"].Range("N105").Value = -11497.6665
$sheets["CRP"].Range("H62").Value = 2978
$sheets["CRP"].Range("I62").Value = 2978
$sheets["CRP"].Range("K62").Value = 2978
$sheets["CRP"].Range("M62").Value = -2354
$sheets["CRP"].Range("H65").Value = 2978
$sheets["CRP"].Range("I65").Value = 2978
$sheets["CRP"].Range("K65").Value = 14890
$sheets["CRP"].Range("M65").Value = -11770
$sheets["CRP"].Range("H107").Value = 1771.091
$sheets["CRP"].Range("I107").Value = 3135.75
$sheets["CRP"].Range("J107").Value = 991.2857
$sheets["CRP"].Range("K107").Value = 3135.75
$sheets["CRP"].Range("L107").Value = 991.2857
$sheets["CRP"].Range("M107").Value = -1215.75
$sheets["CRP"].Range("N107").Value = -4831.2857
$sheets["CUL"].Range("H3").Value = 4459.5625
$sheets["CUL"].Range("I3").Value = 3035.9
$sheets["CUL"].Range("J3").Value = 6832.3335
$sheets["CUL"].Range("K3").Value = 9107.700000000001
$sheets["CUL"].Range("L3").Value = 20497.0005
$sheets["CUL"].Range("M3").Value = -8995.700000000001
$sheets["CUL"].Range("N3").Value = -20721.0005
$sheets["CUL"].Range("H131").Value = 957.34784
$sheets["CUL"].Range("I131").Value = 941.5833
$sheets["CUL"].Range("J131").Value = 974.5454999999999
$sheets["CUL"].Range("K131").Value = 2824.7499
$sheets["CUL"].Range("L131").Value = 2923.6365
$sheets["CUL"].Range("M131").Value = 2215.2501
$sheets["CUL"].Range("N131").Value = -13003.6365
$sheets["GSM"].Range("H33").Value = 10000
$sheets["GSM"].Range("J33").Value = 10000
$sheets["GSM"].Range("L33").Value = 10000
$sheets["GSM"].Range("N33").Value = -10504
$sheets["GSM"].Range("H36").Value = 1249
$sheets["GSM"].Range("I36").Value = 1339
$sheets["GSM"].Range("J36").Value = 979
$sheets["GSM"].Range("K36").Value = 1339
$sheets["GSM"].Range("L36").Value = 979
$sheets["GSM"].Range("M36").Value = -854
$sheets["GSM"].Range("N36").Value = -1949
$sheets["GSM"].Range("H80").Value = 7600
$sheets["GSM"].Range("I80").Value = 9162.5
$sheets["GSM"].Range("J80").Value = 5100
$sheets["GSM"].Range("K80").Value = 9162.5
$sheets["GSM"].Range("L80").Value = 5100
$sheets["GSM"].Range("M80").Value = -8164.5
$sheets["GSM"].Range("N80").Value = -7096
$sheets["GSM"].Range("H83").Value = 7600
$sheets["GSM"].Range("I83").Value = 9162.5
$sheets["GSM"].Range("J83").Value = 5100
$sheets["GSM"].Range("K83").Value = 45812.5
$sheets["GSM"].Range("L83").Value = 25500
$sheets["GSM"].Range("M83").Value = -40820.5
$sheets["GSM"].Range("N83").Value = -35484
$sheets["GSM"].Range("H97").Value = 918.26666
$sheets["GSM"].Range("I97").Value = 922.8333
$sheets["GSM"].Range("J97").Value = 900
$sheets["GSM"].Range("K97").Value = 922.8333
$sheets["GSM"].Range("L97").Value = 900
$sheets["GSM"].Range("M97").Value = -426.8333
$sheets["GSM"].Range("N97").Value = -1892
$sheets["GSM"].Range("H113").Value = 1895.5333
$sheets["GSM"].Range("I113").Value = 1893.909
$sheets["GSM"].Range("J113").Value = 1900
$sheets["GSM"].Range("K113").Value = 1893.909
$sheets["GSM"].Range("L113").Value = 1900
$sheets["GSM"].Range("M113").Value = 276.0909999999999
$sheets["GSM"].Range("N113").Value = -6240
$sheets["LTW"].Range("H16").Value = 1087.6
$sheets["LTW"].Range("I16").Value = 984.5
$sheets["LTW"].Range("J16").Value = 1500
$sheets["LTW"].Range("K16").Value = 984.5
$sheets["LTW"].Range("L16").Value = 1500
$sheets["LTW"].Range("M16").Value = -814.5
$sheets["LTW"].Range("N16").Value = -1840
$sheets["LTW"].Range("H46").Value = 637.75
$sheets["LTW"].Range("I46").Value = 690.5
$sheets["LTW"].Range("J46").Value = 585
$sheets["LTW"].Range("K46").Value = 690.5
$sheets["LTW"].Range("L46").Value = 585
$sheets["LTW"].Range("M46").Value = -502.5
$sheets["LTW"].Range("N46").Value = -961
$sheets["LTW"].Range("H61").Value = 30151.143
$sheets["LTW"].Range("I61").Value = 30151.143
$sheets["LTW"].Range("J61").Value = 0
$sheets["LTW"].Range("K61").Value = 30151.143
$sheets["LTW"].Range("L61").Value = 0
$sheets["LTW"].Range("M61").Value = -29949.143
$sheets["LTW"].Range("N61").ClearContents()
$sheets["LTW"].Range("H82").Value = 853.8182
$sheets["LTW"].Range("I82").Value = 826.55554
$sheets["LTW"].Range("J82").Value = 976.5
$sheets["LTW"].Range("K82").Value = 826.55554
$sheets["LTW"].Range("L82").Value = 976.5
$sheets["LTW"].Range("M82").Value = -465.55554
$sheets["LTW"].Range("N82").Value = -1698.5
$sheets["LTW"].Range("H85").Value = 853.8182
$sheets["LTW"].Range("I85").Value = 826.55554
$sheets["LTW"].Range("J85").Value = 976.5
$sheets["LTW"].Range("K85").Value = 826.55554
$sheets["LTW"].Range("L85").Value = 976.5
$sheets["LTW"].Range("M85").Value = 421.44446
$sheets["LTW"].Range("N85").Value = -3472.5
$sheets["LTW"].Range("H113").Value = 30151.143
$sheets["LTW"].Range("I113").Value = 30151.143
$sheets["LTW"].Range("J113").Value = 0
$sheets["LTW"].Range("K113").Value = 30151.143
$sheets["LTW"].Range("L113").Value = 0
$sheets["LTW"].Range("M113").Value = -27981.143
$sheets["LTW"].Range("N113").ClearContents()
$sheets["WVR"].Range("H62").Value = 3751.4167
$sheets["WVR"].Range("J62").Value = 3502.5
$sheets["WVR"].Range("L62").Value = 3502.5
$sheets["WVR"].Range("N62").Value = -4750.5
$sheets["WVR"].Range("H65").Value = 3751.4167
$sheets["WVR"].Range("J65").Value = 3502.5
$sheets["WVR"].Range("L65").Value = 17512.5
$sheets["WVR"].Range("N65").Value = -23752.5
$sheets["WVR"].Range("H96").Value = 540.4286
$sheets["WVR"].Range("I96").Value = 493.33334
$sheets["WVR"].Range("J96").Value = 575.75
$sheets["WVR"].Range("K96").Value = 493.33334
$sheets["WVR"].Range("L96").Value = 575.75
$sheets["WVR"].Range("M96").Value = 879.66666
$sheets["WVR"].Range("N96").Value = -3321.75
